$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cell values from the crypto price refresh.
# Cells whose new text looks like a plain number (e.g. "575.03") are forced to
# remain text so Excel does not silently convert them to floating point numbers
# (which would lose formatting like trailing zeros or change the stored value).

$ws.Range("D2").Value = "62.827.87"
$ws.Range("E2").Value = "  +0.30%  "
$ws.Range("D3").Value = "2.461.67"
$ws.Range("E3").Value = "  +0.77%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "575.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.76"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.537"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "2.461.04"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.112"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.23%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.162"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.56%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.26"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.82%  "
$ws.Range("E13").Value = "  +0.83%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.91"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.49%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000177"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.54%  "
$ws.Range("D16").Value = "2.908.89"
$ws.Range("E16").Value = "  +0.87%  "
$ws.Range("D17").Value = "62.732.32"
$ws.Range("E17").Value = "  +0.34%  "
$ws.Range("D18").Value = "2.461.05"
$ws.Range("E18").Value = "  +0.94%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.05"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.98"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.77%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "326.86"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.29%  "
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("E23").Value = "  +9.64%  "
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.25"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +20.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "65.70"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.68%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "654.02"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.06%  "
$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").Value = "0.0₃0978"
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("B29").Value = "WrappedeETH"
$ws.Range("C29").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D29").Value = "2.583.04"
$ws.Range("E29").Value = "  +0.95%  "
$ws.Range("E30").Value = "  -12.23%  "
$ws.Range("E31").Value = "  +2.65%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.19%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.85"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.92%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.134"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.40%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.53"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.76"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.38%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.370"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.79%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.71"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.84%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "150.89"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.36"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.01%  "
$ws.Range("E42").Value = "  +1.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.73"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.50%  "
$ws.Range("D44").Value = "0.0₆0313"
$ws.Range("E44").Value = "  -52.94%  "
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "154.38"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +7.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "15.23"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.45%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.59"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.17%  "
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "20.31"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.01%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.606"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.41%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0510"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.28%  "
